$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new log entry (row 3 of Table1)
$ws.Range("A3").Value = 45909
$ws.Range("B3").Value = 0.41666666666666669
$ws.Range("C3").Value = 0.58333333333333337
$ws.Range("E3").Value = "Finished up Chapter 3 from the book, then did Chapter 4. This took a bit longer than expected, the chatpers were big. Also installing packages for later took some time. Chapter 4 included the first practice program. Figuring this out took a bit since it is the first program I have written in Python. The biggest challenge here was using Python correctly, not necessarily the logic behind the question. In the end I made a working program."

# Row grows tall to show the wrapped description text (matches Excel's
# automatic row-height recalculation after typing a long wrapped entry)
$ws.Rows.Item(3).RowHeight = 72

# Move the active selection to the next empty row, as it was left after entry
[void]$ws.Range("E4").Select()
